# Big update with local data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schedule dates for lecture rows 8 and 9 (Table13).
# Dependent formula columns (Day = Date, HW Deadline = Date+7) recalc automatically.
$ws.Range("C8").Value = 42437
$ws.Range("C9").Value = 42438

# Match the author's final selection/active cell in the sheet.
[void]$ws.Range("C9").Select()
